# Update auto-increase logic projection results to reflect the
# "new hire only" vs "all eligible" auto-increase scenario recalculation.
# Underlying eligibility/participation counts (and all derived metrics)
# shift for projection years 2020-2024 (rows 2-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1065069431212171
$ws.Range("H2").Value = 0.09084415736809691
$ws.Range("I2").Value = 528276.0912988938
$ws.Range("J2").Value = 201524.0458504469
$ws.Range("L2").Value = 201524.0458504469
$ws.Range("M2").Value = 729800.1371493406
$ws.Range("N2").Value = 10161399.9688
$ws.Range("O2").Value = 9753659.038699998
$ws.Range("P2").Value = 0.01983231114504054
$ws.Range("Q2").Value = 0.02066137898104204
$ws.Range("C3").Value = 103
$ws.Range("D3").Value = 88
$ws.Range("E3").Value = 0.8543689320388349
$ws.Range("F3").Value = 0.8543689320388349
$ws.Range("G3").Value = 0.1101170774262103
$ws.Range("H3").Value = 0.094080609839869
$ws.Range("I3").Value = 595699.0697549444
$ws.Range("J3").Value = 232637.1150845023
$ws.Range("L3").Value = 232637.1150845023
$ws.Range("M3").Value = 828336.1848394468
$ws.Range("N3").Value = 10536239.823564
$ws.Range("O3").Value = 10128866.665561
$ws.Range("P3").Value = 0.02207970955294848
$ws.Range("Q3").Value = 0.0229677339791122
$ws.Range("C4").Value = 103
$ws.Range("D4").Value = 87
$ws.Range("E4").Value = 0.8446601941747572
$ws.Range("F4").Value = 0.8365384615384616
$ws.Range("G4").Value = 0.1155345313075407
$ws.Range("H4").Value = 0.09664907907457733
$ws.Range("I4").Value = 659315.5951792673
$ws.Range("J4").Value = 257359.0052028746
$ws.Range("L4").Value = 257359.0052028746
$ws.Range("M4").Value = 916674.600382142
$ws.Range("N4").Value = 10860875.61657092
$ws.Range("O4").Value = 10452531.26382783
$ws.Range("P4").Value = 0.02369597206418703
$ws.Range("Q4").Value = 0.02462169198129975
$ws.Range("D5").Value = 85
$ws.Range("E5").Value = 0.8095238095238095
$ws.Range("F5").Value = 0.8095238095238095
$ws.Range("G5").Value = 0.1172245253426888
$ws.Range("H5").Value = 0.09489604432503375
$ws.Range("I5").Value = 678472.924129093
$ws.Range("J5").Value = 264363.7059061846
$ws.Range("L5").Value = 264363.7059061846
$ws.Range("M5").Value = 942836.6300352775
$ws.Range("N5").Value = 11360984.84936805
$ws.Range("O5").Value = 10950290.16604267
$ws.Range("P5").Value = 0.02326943565292139
$ws.Range("Q5").Value = 0.02414216444473664
$ws.Range("C6").Value = 105
$ws.Range("E6").Value = 0.8571428571428571
$ws.Range("G6").Value = 0.1134076443070747
$ws.Range("H6").Value = 0.09628950931732752
$ws.Range("I6").Value = 716159.0169867697
$ws.Range("J6").Value = 280145.5696502721
$ws.Range("L6").Value = 280145.5696502721
$ws.Range("M6").Value = 996304.5866370418
$ws.Range("N6").Value = 11698901.34244909
$ws.Range("O6").Value = 11284435.81862395
$ws.Range("P6").Value = 0.02394631439738472
$ws.Range("Q6").Value = 0.02482583747677638
